# Update R script paths for PBA40 -> RTP plus IPA model run names
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTP2021")

# Row 2: was 2035 / 2035_TM151_IPA_loPop_loAOC_00 / IP
$ws.Range("A2").Value = 2005
$ws.Range("B2").Value = "2005_TM152_IPA_02"
$ws.Range("C2").Value = "IP"

# Row 3: was 2035_TM151_IPA_loPop_hiAOC_00 / IP_hiAOC
$ws.Range("B3").Value = "2035_TM152_IPA_aoc1421_00"
$ws.Range("C3").Value = "IP_aoc1421"

# Row 4: was 2035_TM151_IPA_hiPop_loAOC_00 / IP_hiPop
$ws.Range("B4").Value = "2035_TM152_IPA_aoc1562_00"
$ws.Range("C4").Value = "IP_aoc1562"

# Row 5: was 2035_TM151_IPA_hiPop_hiAOC_00 / IP_hiPop_hiAOC
$ws.Range("B5").Value = "2035_TM152_IPA_aoc1795_00"
$ws.Range("C5").Value = "IP_aoc1795"

# Move the active selection from C4 to C5
$ws.Range("C5").Select()
